# Apply the edits described by the diff:
#  - Rows 2-7 get new Question (col A) values, and col B is set to the AI text
#    for every row (2-7).
#  - Rows 8-15 are removed entirely (the sheet now only spans A1:B7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aiText = "AI adalah Artificial Intelligence, yaitu kecerdasan buatan yang memungkinkan mesin untuk meniru kecerdasan manusia dalam menyelesaikan tugas."

# Delete rows 8 through 15 first (from bottom up) so row numbers of earlier
# rows are unaffected while removing them.
$ws.Range("A8:B15").EntireRow.Delete()

# Update the remaining question/answer rows (2-7).
$ws.Range("A2").Value = "go"
$ws.Range("B2").Value = $aiText

$ws.Range("A3").Value = "ai"
$ws.Range("B3").Value = $aiText

$ws.Range("A4").Value = "acid"
$ws.Range("B4").Value = $aiText

$ws.Range("A5").Value = "apa itu acid"
$ws.Range("B5").Value = $aiText

$ws.Range("A6").Value = "apa itu ai"
$ws.Range("B6").Value = $aiText

$ws.Range("A7").Value = "exit "
$ws.Range("B7").Value = $aiText
